$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("teppppppp", 22, "Subscribed", "Employed"),
    @("asdf", 33, "Subscribed", "Unemployed"),
    @("patrick", 44, "Not Subscribed", "Unemployed"),
    @("asdfasdfasdfasdfjklhasdfkjlashdf", 44, "Subscribed", "Unemployed")
)

$startRow = 56
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
